# Apply the "first full GUI version" update to paths.xlsx
# - D2 used to hold the number 1; it now holds the text "0"
#   (mirrors the is_random column switching from a numeric flag to the
#   GUI's string output).
# - The active selection moves from C7 to D2, reflecting where the user
#   was last working when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Store "0" as text (not the number 0) in D2.
$ws.Range("D2").Value = "0"

# Update the saved selection/active cell to D2.
$ws.Range("D2").Select()
